$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.54690179691499
$ws.Range("C2").Value = 9.428555705286524
$ws.Range("D2").Value = 3.966546793897766
$ws.Range("E2").Value = 11.535794073821
$ws.Range("F2").Value = 20.30730892550126
$ws.Range("I2").Value = 17.13428314545414
$ws.Range("M2").Value = 14.54929968178754
$ws.Range("O2").Value = 18.06989336467801
$ws.Range("B3").Value = 11.86650553936523
$ws.Range("C3").Value = 8.92275859314053
$ws.Range("D3").Value = 3.910455146537271
$ws.Range("E3").Value = 11.47162278415189
$ws.Range("F3").Value = 20.31224683771581
$ws.Range("I3").Value = 17.27335154175147
$ws.Range("M3").Value = 14.22167374363967
$ws.Range("O3").Value = 18.14189883606086
$ws.Range("B4").Value = 11.4279746410797
$ws.Range("C4").Value = 8.595601600096947
$ws.Range("D4").Value = 3.875358465795599
$ws.Range("E4").Value = 11.43735208872515
$ws.Range("F4").Value = 20.32376456729344
$ws.Range("I4").Value = 17.36395814447198
$ws.Range("M4").Value = 14.01934933536986
$ws.Range("O4").Value = 18.19256791516396
$ws.Range("B5").Value = 11.24419623374692
$ws.Range("C5").Value = 8.458185047238223
$ws.Range("D5").Value = 3.860901676604927
$ws.Range("E5").Value = 11.42468621781445
$ws.Range("F5").Value = 20.33058498205356
$ws.Range("I5").Value = 17.40219134464362
$ws.Range("M5").Value = 13.93672711513706
$ws.Range("O5").Value = 18.21482996576233
$ws.Range("B6").Value = 11.21337845651533
$ws.Range("C6").Value = 8.435122191506844
$ws.Range("D6").Value = 3.858492108761238
$ws.Range("E6").Value = 11.42266178314034
$ws.Range("F6").Value = 20.33184571673693
$ws.Range("I6").Value = 17.40861902445322
$ws.Range("M6").Value = 13.92300082452816
$ws.Range("O6").Value = 18.21862374933642
$ws.Range("B7").Value = 11.42551646631564
$ws.Range("C7").Value = 8.593764830893246
$ws.Range("D7").Value = 3.87516410856585
$ws.Range("E7").Value = 11.43717599934872
$ws.Range("F7").Value = 20.32384794913849
$ws.Range("I7").Value = 17.36446846790269
$ws.Range("M7").Value = 14.0182356058303
$ws.Range("O7").Value = 18.1928616272523
$ws.Range("B8").Value = 12.31669087702584
$ws.Range("C8").Value = 9.257641844380462
$ws.Range("D8").Value = 3.947348735587962
$ws.Range("E8").Value = 11.5126095161764
$ws.Range("F8").Value = 20.30724673972433
$ws.Range("I8").Value = 17.18114978252272
$ws.Range("M8").Value = 14.43664501186319
$ws.Range("O8").Value = 18.09337485398054
$ws.Range("B9").Value = 13.89437944252291
$ws.Range("C9").Value = 10.42524318532942
$ws.Range("D9").Value = 4.083254698881348
$ws.Range("E9").Value = 11.70068904893028
$ws.Range("F9").Value = 20.34223137451639
$ws.Range("I9").Value = 16.86315111280821
$ws.Range("M9").Value = 15.24294947752668
$ws.Range("O9").Value = 17.94993617654536
$ws.Range("B10").Value = 14.94460345886194
$ws.Range("C10").Value = 11.19862334985475
$ws.Range("D10").Value = 4.179094439902826
$ws.Range("E10").Value = 11.86242023282502
$ws.Range("F10").Value = 20.40923989164077
$ws.Range("I10").Value = 16.65493346427622
$ws.Range("M10").Value = 15.82037789357161
$ws.Range("O10").Value = 17.87658233447664
$ws.Range("B11").Value = 15.39796299299143
$ws.Range("C11").Value = 11.53175679950014
$ws.Range("D11").Value = 4.221710528619226
$ws.Range("E11").Value = 11.94085970524838
$ws.Range("F11").Value = 20.44867126401439
$ws.Range("I11").Value = 16.56575881297928
$ws.Range("M11").Value = 16.07859846550222
$ws.Range("O11").Value = 17.85027171174593
$ws.Range("B12").Value = 15.56608437222395
$ws.Range("C12").Value = 11.65519864166581
$ws.Range("D12").Value = 4.237697966628862
$ws.Range("E12").Value = 11.97123840520256
$ws.Range("F12").Value = 20.46488404677566
$ws.Range("I12").Value = 16.53279100374435
$ws.Range("M12").Value = 16.17564655711247
$ws.Range("O12").Value = 17.84133124336466
$ws.Range("B13").Value = 15.53003523569674
$ws.Range("C13").Value = 11.62873398901254
$ws.Range("D13").Value = 4.234261613931348
$ws.Range("E13").Value = 11.96466618142935
$ws.Range("F13").Value = 20.46133548112097
$ws.Range("I13").Value = 16.539855543961
$ws.Range("M13").Value = 16.15477960181212
$ws.Range("O13").Value = 17.8432111268336
$ws.Range("B14").Value = 15.41186594663047
$ws.Range("C14").Value = 11.5419668245271
$ws.Range("D14").Value = 4.223028888942523
$ws.Range("E14").Value = 11.94334559077706
$ws.Range("F14").Value = 20.44997946887271
$ws.Range("I14").Value = 16.56303046417892
$ws.Range("M14").Value = 16.08659790346296
$ws.Range("O14").Value = 17.84951563333689
$ws.Range("B15").Value = 15.33901949575852
$ws.Range("C15").Value = 11.48846613575926
$ws.Range("D15").Value = 4.216128676032812
$ws.Range("E15").Value = 11.93037330717211
$ws.Range("F15").Value = 20.44319018936567
$ws.Range("I15").Value = 16.57733015036195
$ws.Range("M15").Value = 16.04473632765114
$ws.Range("O15").Value = 17.85351074405762
$ws.Range("B16").Value = 14.91448010943391
$ws.Range("C16").Value = 11.17647434534285
$ws.Range("D16").Value = 4.176288827940787
$ws.Range("E16").Value = 11.8573899165458
$ws.Range("F16").Value = 20.40684258743209
$ws.Range("I16").Value = 16.66087309395177
$ws.Range("M16").Value = 15.80340475758275
$ws.Range("O16").Value = 17.87844451041102
$ws.Range("B17").Value = 14.64775331684458
$ws.Range("C17").Value = 10.98027504190365
$ws.Range("D17").Value = 4.15159023928968
$ws.Range("E17").Value = 11.81384674044089
$ws.Range("F17").Value = 20.38683313349658
$ws.Range("I17").Value = 16.71354643489199
$ws.Range("M17").Value = 15.6541458291898
$ws.Range("O17").Value = 17.89555415004024
$ws.Range("B18").Value = 14.49204526988545
$ws.Range("C18").Value = 10.86566894836814
$ws.Range("D18").Value = 4.137292365557045
$ws.Range("E18").Value = 11.7892615742552
$ws.Range("F18").Value = 20.3761670396912
$ws.Range("I18").Value = 16.74436441596619
$ws.Range("M18").Value = 15.56788195363548
$ws.Range("O18").Value = 17.90605910556112
$ws.Range("B19").Value = 14.43893289096704
$ws.Range("C19").Value = 10.82656422903277
$ws.Range("D19").Value = 4.132435853305035
$ws.Range("E19").Value = 11.78101713145296
$ws.Range("F19").Value = 20.37270057656994
$ws.Range("I19").Value = 16.75488835209151
$ws.Range("M19").Value = 15.53860633787331
$ws.Range("O19").Value = 17.90972968774737
$ws.Range("B20").Value = 14.67638463013784
$ws.Range("C20").Value = 11.00134286149211
$ws.Range("D20").Value = 4.154229025037042
$ws.Range("E20").Value = 11.81843457838207
$ws.Range("F20").Value = 20.38887598053883
$ws.Range("I20").Value = 16.70788525317048
$ws.Range("M20").Value = 15.67007825476686
$ws.Range("O20").Value = 17.89366402545276
$ws.Range("B21").Value = 15.44667197253522
$ws.Range("C21").Value = 11.56752610205544
$ws.Range("D21").Value = 4.226332363879147
$ws.Range("E21").Value = 11.94958984030691
$ws.Range("F21").Value = 20.4532802982669
$ws.Range("I21").Value = 16.5562016696355
$ws.Range("M21").Value = 16.10664516828393
$ws.Range("O21").Value = 17.84763602965585
$ws.Range("B22").Value = 15.92935687388675
$ws.Range("C22").Value = 11.92176558412559
$ws.Range("D22").Value = 4.272575763138382
$ws.Range("E22").Value = 12.03923133753445
$ws.Range("F22").Value = 20.50283434154344
$ws.Range("I22").Value = 16.4617364704722
$ws.Range("M22").Value = 16.38764722720544
$ws.Range("O22").Value = 17.82351885165479
$ws.Range("B23").Value = 15.67365003035578
$ws.Range("C23").Value = 11.73415234308487
$ws.Range("D23").Value = 4.247978253241791
$ws.Range("E23").Value = 11.99103751452025
$ws.Range("F23").Value = 20.47570619317558
$ws.Range("I23").Value = 16.51172595601765
$ws.Range("M23").Value = 16.23809566929809
$ws.Range("O23").Value = 17.83584245701244
$ws.Range("B24").Value = 14.66344777067229
$ws.Range("C24").Value = 10.99182373166206
$ws.Range("D24").Value = 4.153036335554735
$ws.Range("E24").Value = 11.81635901832842
$ws.Range("F24").Value = 20.38794979985137
$ws.Range("I24").Value = 16.71044300595869
$ws.Range("M24").Value = 15.66287660872992
$ws.Range("O24").Value = 17.89451646932068
$ws.Range("B25").Value = 13.48638658319317
$ws.Range("C25").Value = 10.12405879918552
$ws.Range("D25").Value = 4.047153096413001
$ws.Range("E25").Value = 11.64559573933182
$ws.Range("F25").Value = 20.32551366280235
$ws.Range("I25").Value = 16.94472447365861
$ws.Range("M25").Value = 15.02701620498745
$ws.Range("O25").Value = 17.98315233507081
